$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "['LWSA','QUAL','SUZB']"
$ws.Range("C4").Value = "['EMBR','BRFS','ELET']"
$ws.Range("C5").Value = "['RADL','GOLL','VALE']"
$ws.Range("C6").Value = "['BPAC','LWSA','QUAL']"
$ws.Range("C7").Value = "['JBSS','UGPA','BRKM']"

$ws.Range("C8").Select()
